$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2153846153846154
$ws.Range("C2").Value = 0.5115384615384615
$ws.Range("J2").Value = 0.02692307692307692
$ws.Range("P2").Value = 0.1384615384615385
$ws.Range("S2").Value = 0.1076923076923077
$ws.Range("B3").Value = 0.007407407407407408
$ws.Range("C3").Value = 0.02962962962962963
$ws.Range("J3").Value = 0.02962962962962963
$ws.Range("P3").Value = 0.7407407407407407
$ws.Range("S3").Value = 0.1925925925925926
$ws.Range("J4").Value = 0.0784313725490196
$ws.Range("P4").Value = 0.5882352941176471
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.8
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.03703703703703703
$ws.Range("D6").Value = 0.01587301587301587
$ws.Range("E6").Value = 0.005291005291005291
$ws.Range("F6").Value = 0.04232804232804233
$ws.Range("J6").Value = 0.2063492063492063
$ws.Range("O6").Value = 0.02116402116402116
$ws.Range("Q6").Value = 0.2433862433862434
$ws.Range("R6").Value = 0.06878306878306878
$ws.Range("S6").Value = 0.3597883597883598
$ws.Range("B7").Value = 0.06666666666666667
$ws.Range("D7").Value = 0.01904761904761905
$ws.Range("E7").Value = 0.004761904761904762
$ws.Range("F7").Value = 0.06666666666666667
$ws.Range("J7").Value = 0.1
$ws.Range("O7").Value = 0.009523809523809525
$ws.Range("Q7").Value = 0.219047619047619
$ws.Range("R7").Value = 0.1
$ws.Range("S7").Value = 0.4142857142857143
$ws.Range("B8").Value = 0.09221311475409837
$ws.Range("D8").Value = 0.01844262295081967
$ws.Range("F8").Value = 0.04918032786885246
$ws.Range("J8").Value = 0.09221311475409837
$ws.Range("O8").Value = 0.01639344262295082
$ws.Range("Q8").Value = 0.2110655737704918
$ws.Range("R8").Value = 0.07581967213114754
$ws.Range("S8").Value = 0.444672131147541
$ws.Range("B9").Value = 0.1118421052631579
$ws.Range("D9").Value = 0.006578947368421052
$ws.Range("E9").Value = 0.006578947368421052
$ws.Range("F9").Value = 0.05921052631578947
$ws.Range("J9").Value = 0.09868421052631579
$ws.Range("O9").Value = 0.02631578947368421
$ws.Range("Q9").Value = 0.1710526315789474
$ws.Range("R9").Value = 0.06578947368421052
$ws.Range("S9").Value = 0.4539473684210527
$ws.Range("B10").Value = 0.1051212938005391
$ws.Range("D10").Value = 0.03054806828391734
$ws.Range("E10").Value = 0.001796945193171608
$ws.Range("F10").Value = 0.07367475292003593
$ws.Range("J10").Value = 0.1132075471698113
$ws.Range("O10").Value = 0.01796945193171608
$ws.Range("Q10").Value = 0.2309074573225517
$ws.Range("R10").Value = 0.06918238993710692
$ws.Range("S10").Value = 0.3575920934411501
$ws.Range("F11").Value = 0.003105590062111801
$ws.Range("G11").Value = 0.1459627329192547
$ws.Range("J11").Value = 0.08695652173913043
$ws.Range("K11").Value = 0.1925465838509317
$ws.Range("L11").Value = 0.5527950310559007
$ws.Range("S11").Value = 0.01863354037267081
$ws.Range("G12").Value = 0.770949720670391
$ws.Range("J12").Value = 0.1955307262569832
$ws.Range("K12").Value = 0.0111731843575419
$ws.Range("L12").Value = 0.0111731843575419
$ws.Range("S12").Value = 0.0111731843575419
$ws.Range("F13").Value = 0.02272727272727273
$ws.Range("G13").Value = 0.7045454545454546
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.02272727272727273
$ws.Range("F15").Value = 0.01030927835051546
$ws.Range("H15").Value = 0.211340206185567
$ws.Range("I15").Value = 0.06701030927835051
$ws.Range("J15").Value = 0.3505154639175257
$ws.Range("K15").Value = 0.05670103092783505
$ws.Range("M15").Value = 0.01030927835051546
$ws.Range("O15").Value = 0.03608247422680412
$ws.Range("S15").Value = 0.2577319587628866
$ws.Range("F16").Value = 0.02439024390243903
$ws.Range("H16").Value = 0.1951219512195122
$ws.Range("I16").Value = 0.1036585365853658
$ws.Range("J16").Value = 0.4024390243902439
$ws.Range("K16").Value = 0.1280487804878049
$ws.Range("M16").Value = 0.01829268292682927
$ws.Range("O16").Value = 0.04878048780487805
$ws.Range("S16").Value = 0.07926829268292683
$ws.Range("F17").Value = 0.01470588235294118
$ws.Range("H17").Value = 0.203781512605042
$ws.Range("I17").Value = 0.07983193277310924
$ws.Range("J17").Value = 0.4138655462184874
$ws.Range("K17").Value = 0.1092436974789916
$ws.Range("M17").Value = 0.02100840336134454
$ws.Range("O17").Value = 0.06302521008403361
$ws.Range("S17").Value = 0.09453781512605042
$ws.Range("F18").Value = 0.03225806451612903
$ws.Range("H18").Value = 0.2258064516129032
$ws.Range("I18").Value = 0.08387096774193549
$ws.Range("J18").Value = 0.3870967741935484
$ws.Range("K18").Value = 0.1483870967741935
$ws.Range("M18").Value = 0.01290322580645161
$ws.Range("O18").Value = 0.01290322580645161
$ws.Range("S18").Value = 0.0967741935483871
$ws.Range("F19").Value = 0.00782608695652174
$ws.Range("H19").Value = 0.2504347826086957
$ws.Range("I19").Value = 0.06260869565217392
$ws.Range("J19").Value = 0.3539130434782609
$ws.Range("K19").Value = 0.1252173913043478
$ws.Range("M19").Value = 0.02434782608695652
$ws.Range("O19").Value = 0.07130434782608695
$ws.Range("S19").Value = 0.1043478260869565
